# Doing Updates for Financials
# Insert a new "most recent period" column (D) in front of the existing
# yearly data, shifting the old D:K columns to E:L, then populate the new
# column with the latest period's figures. A handful of rows also get a
# correction to the (now-shifted) E column where the historical figure
# itself changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new column before D; this shifts D:K -> E:L (values, styles,
# and types all move together, same as Excel's own Insert command).
$ws.Columns("D:D").Insert()

# New column D values, keyed by row.
$newD = @{
    7   = 43465
    8   = 1065500
    9   = 877100
    10  = 188400
    12  = 19200
    14  = 55100
    15  = 4000
    17  = 1053900
    18  = 11700
    20  = 30400
    21  = 75900
    22  = 5700
    23  = 36400
    24  = 11500
    26  = 24800
    27  = 24800
    29  = 0
    32  = -30400
    33  = 24800
    35  = 24800
    38  = 43465
    41  = 34400
    80  = 43465
    43  = 131500
    44  = 93800
    45  = 9600
    46  = 269300
    47  = 84600
    48  = 228400
    49  = 117700
    52  = 7400
    54  = 707400
    57  = 112800
    59  = 42500
    60  = 155300
    61  = 101500
    62  = 95800
    66  = 352500
    72  = 497500
    76  = 354900
    81  = 24800
    83  = 33800
    89  = 97800
    91  = -40800
    94  = -34100
    96  = -14600
    100 = -64100
    101 = -1700
    102 = -2100
}

foreach ($row in $newD.Keys) {
    $ws.Range("D$row").Value = $newD[$row]
}

# Rows 7, 38 and 80 hold "Period Ending" dates - format the new column D
# like the adjacent (shifted) date cells.
$ws.Range("D7").NumberFormat = $ws.Range("E7").NumberFormat
$ws.Range("D38").NumberFormat = $ws.Range("E38").NumberFormat
$ws.Range("D80").NumberFormat = $ws.Range("E80").NumberFormat

# A few rows have a genuine correction to the historical figure that now
# sits in column E (not merely the value that shifted over from the old
# column D), so overwrite those explicitly.
$ws.Range("E43").Value = 162100
$ws.Range("E48").Value = 994000
$ws.Range("E49").Value = 209300
$ws.Range("E59").Value = 84900
